# Apply updated crypto price/volume data to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.253.76"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "1.646.66"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  -0.09%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "217.09"
$r.ClearFormats()
$ws.Range("E5").Value = "  +0.41%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.506"
$r.ClearFormats()
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.07%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.0637"
$r.ClearFormats()
$ws.Range("E9").Value = "  +0.14%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "19.90"
$r.ClearFormats()
$ws.Range("E10").Value = "  +0.88%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0795"
$r.ClearFormats()
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "4.30"
$r.ClearFormats()
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "1.874.55"
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "1.633.28"
$ws.Range("E14").Value = "  -0.40%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.548"
$r.ClearFormats()
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").Value = "0.0₃0766"
$ws.Range("E16").Value = "  -0.21%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "63.33"
$r.ClearFormats()
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "26.246.57"
$ws.Range("E18").Value = "  +1.52%  "
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "4.44"
$r.ClearFormats()
$ws.Range("E20").Value = "  -0.85%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "195.36"
$r.ClearFormats()
$ws.Range("E21").Value = "  +1.32%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "10.06"
$r.ClearFormats()
$ws.Range("E22").Value = "  +0.80%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "6.32"
$r.ClearFormats()
$ws.Range("E23").Value = "  -0.58%  "
$ws.Range("E24").Value = "  -1.94%  "
$ws.Range("B25").Value = "BinanceUSD"
$ws.Range("C25").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.ClearFormats()
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "142.93"
$r.ClearFormats()
$ws.Range("E26").Value = "  +0.61%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "0.125"
$r.ClearFormats()
$ws.Range("E27").Value = "  +1.23%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "6.98"
$r.ClearFormats()
$ws.Range("E28").Value = "  +0.43%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "15.62"
$r.ClearFormats()
$ws.Range("E29").Value = "  +0.58%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "1.25"
$r.ClearFormats()
$ws.Range("E30").Value = "  +1.12%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.0504"
$r.ClearFormats()
$ws.Range("E31").Value = "  +2.14%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.35"
$r.ClearFormats()
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("E35").Value = "  +1.26%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.913"
$r.ClearFormats()
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("B37").Value = "Maker"
$ws.Range("C37").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D37").Value = "1.134.99"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.555"
$r.ClearFormats()
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("E39").Value = "  -1.14%  "
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("E41").Value = "  +0.05%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "5.55"
$r.ClearFormats()
$ws.Range("E42").Value = "  -0.83%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "100.38"
$r.ClearFormats()
$ws.Range("E43").Value = "  -0.27%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.802"
$r.ClearFormats()
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("D45").Value = "1.784.49"
$ws.Range("E45").Value = "  +0.52%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "56.38"
$r.ClearFormats()
$ws.Range("E46").Value = "  +1.74%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "1.49"
$r.ClearFormats()
$ws.Range("E47").Value = "  +2.95%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.0518"
$r.ClearFormats()
$ws.Range("E48").Value = "  +3.12%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "7.71"
$r.ClearFormats()
$ws.Range("E49").Value = "  +2.98%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.418"
$r.ClearFormats()
$ws.Range("E50").Value = "  +0.25%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.0968"
$r.ClearFormats()
$ws.Range("E51").Value = "  +1.78%  "
